# Refresh market-price-derived profit figures across the Leve profit sheets.
# Values mirror an upstream market-board data refresh (scheduled runner).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 489.55554
$ws.Range("I2").Value = 489.55554
$ws.Range("K2").Value = 489.55554
$ws.Range("M2").Value = -376.55554

# Row 17
$ws.Range("H17").Value = 875336.6
$ws.Range("J17").Value = 928862.2
$ws.Range("L17").Value = 2786586.6
$ws.Range("N17").Value = -2786922.6

# Row 51
$ws.Range("H51").Value = 15175.637
$ws.Range("I51").Value = 18117.125
$ws.Range("J51").Value = 7331.6665
$ws.Range("K51").Value = 18117.125
$ws.Range("L51").Value = 7331.6665
$ws.Range("M51").Value = -17633.125
$ws.Range("N51").Value = -8299.666499999999

# Row 62
$ws.Range("H62").Value = 28581.46
$ws.Range("I62").Value = 3593.3333
$ws.Range("K62").Value = 3593.3333
$ws.Range("M62").Value = -2969.3333

# Row 65
$ws.Range("H65").Value = 28581.46
$ws.Range("I65").Value = 3593.3333
$ws.Range("K65").Value = 17966.6665
$ws.Range("M65").Value = -14846.6665

# Row 100
$ws.Range("H100").Value = 8116842.5
$ws.Range("J100").Value = 26768.826
$ws.Range("L100").Value = 26768.826
$ws.Range("N100").Value = -27850.826

# Row 116
$ws.Range("H116").Value = 11390928
$ws.Range("I116").Value = 13920246
$ws.Range("K116").Value = 13920246
$ws.Range("M116").Value = -13916804

# Row 132
$ws.Range("H132").Value = 2754.1746
$ws.Range("I132").Value = 2569.1897
$ws.Range("K132").Value = 7707.5691
$ws.Range("M132").Value = -5177.5691

# Row 135
$ws.Range("H135").Value = 1929.4595
$ws.Range("I135").Value = 2012.2354
$ws.Range("K135").Value = 18110.1186
$ws.Range("M135").Value = -15575.1186

# Row 137
$ws.Range("H137").Value = 10634.667
$ws.Range("I137").Value = 13561.895
$ws.Range("J137").Value = 3682.5
$ws.Range("K137").Value = 40685.685
$ws.Range("L137").Value = 11047.5
$ws.Range("M137").Value = -38135.685
$ws.Range("N137").Value = -16147.5

# Row 138
$ws.Range("H138").Value = 261494.7
$ws.Range("J138").Value = 4626.483
$ws.Range("L138").Value = 13879.449
$ws.Range("N138").Value = -24159.449

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 6228.0635
$ws.Range("I32").Value = 6315.984
$ws.Range("J32").Value = 777
$ws.Range("K32").Value = 6315.984
$ws.Range("L32").Value = 777
$ws.Range("M32").Value = -6028.984
$ws.Range("N32").Value = -1351

# Row 61
$ws.Range("H61").Value = 5112.0977
$ws.Range("I61").Value = 5108.0356
$ws.Range("K61").Value = 5108.0356
$ws.Range("M61").Value = -4896.0356

# Row 74
$ws.Range("H74").Value = 4071.6128
$ws.Range("I74").Value = 2183.5217
$ws.Range("K74").Value = 2183.5217
$ws.Range("M74").Value = -1309.5217

# Row 77
$ws.Range("H77").Value = 4071.6128
$ws.Range("I77").Value = 2183.5217
$ws.Range("K77").Value = 10917.6085
$ws.Range("M77").Value = -6549.608499999998

# Row 105
$ws.Range("H105").Value = 167375.25
$ws.Range("J105").Value = 167375.25
$ws.Range("L105").Value = 167375.25
$ws.Range("N105").Value = -174363.25

# Row 110
$ws.Range("H110").Value = 3726.0588
$ws.Range("I110").Value = 3161.5557
$ws.Range("K110").Value = 3161.5557
$ws.Range("M110").Value = -1116.5557

# Row 133
$ws.Range("H133").Value = 67999
$ws.Range("I133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("M133").ClearContents()

# Row 134
$ws.Range("H134").Value = 397777.2
$ws.Range("J134").Value = 397777.2
$ws.Range("L134").Value = 397777.2
$ws.Range("N134").Value = -407917.2

# Row 136
$ws.Range("H136").Value = 5112.0977
$ws.Range("I136").Value = 5108.0356
$ws.Range("K136").Value = 15324.1068
$ws.Range("M136").Value = -12774.1068

$ws = $wb.Worksheets.Item("BSM")
# Row 107
$ws.Range("H107").Value = 1635.238
$ws.Range("I107").Value = 1533.7368
$ws.Range("J107").Value = 2599.5
$ws.Range("K107").Value = 1533.7368
$ws.Range("L107").Value = 2599.5
$ws.Range("M107").Value = 386.2632000000001
$ws.Range("N107").Value = -6439.5

# Row 138
$ws.Range("H138").Value = 142439.55
$ws.Range("J138").Value = 142439.55
$ws.Range("L138").Value = 142439.55
$ws.Range("N138").Value = -152719.55

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 3269.8044
$ws.Range("I31").Value = 2823.7
$ws.Range("K31").Value = 2823.7
$ws.Range("M31").Value = -2528.7

# Row 34
$ws.Range("H34").Value = 3269.8044
$ws.Range("I34").Value = 2823.7
$ws.Range("K34").Value = 2823.7
$ws.Range("M34").Value = -2621.7

# Row 35
$ws.Range("H35").Value = 2084.7
$ws.Range("I35").Value = 2084.7
$ws.Range("K35").Value = 2084.7
$ws.Range("M35").Value = -1790.7

# Row 58
$ws.Range("H58").Value = 2006.25
$ws.Range("I58").Value = 1707.9375
$ws.Range("K58").Value = 1707.9375
$ws.Range("M58").Value = -1504.9375

# Row 122
$ws.Range("H122").Value = 5512.1
$ws.Range("J122").Value = 1763.4546
$ws.Range("L122").Value = 5290.3638
$ws.Range("N122").Value = -10190.3638

# Row 132
$ws.Range("H132").Value = 14380.625
$ws.Range("I132").Value = 1362.48
$ws.Range("J132").Value = 60874
$ws.Range("K132").Value = 4087.44
$ws.Range("L132").Value = 182622
$ws.Range("M132").Value = -1557.44
$ws.Range("N132").Value = -187682

# Row 136
$ws.Range("H136").Value = 2006.25
$ws.Range("I136").Value = 1707.9375
$ws.Range("K136").Value = 5123.8125
$ws.Range("M136").Value = -2573.8125

$ws = $wb.Worksheets.Item("CUL")
# Row 22
$ws.Range("H22").Value = 990
$ws.Range("I22").Value = 990
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 2970
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()
$ws.Range("M22").Value = -2801

# Row 25
$ws.Range("H25").Value = 2300.2942
$ws.Range("J25").Value = 2381.5625
$ws.Range("L25").Value = 7144.6875
$ws.Range("N25").Value = -7482.6875

# Row 26
$ws.Range("H26").Value = 62.916668
$ws.Range("I26").Value = 61.666668
$ws.Range("K26").Value = 185.000004
$ws.Range("M26").Value = 102.999996

# Row 27
$ws.Range("H27").Value = 990
$ws.Range("I27").Value = 990
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 2970
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()
$ws.Range("M27").Value = -2868

# Row 30
$ws.Range("H30").Value = 2300.2942
$ws.Range("J30").Value = 2381.5625
$ws.Range("L30").Value = 7144.6875
$ws.Range("N30").Value = -7348.6875

# Row 54
$ws.Range("H54").Value = 10749.75
$ws.Range("J54").Value = 13999.667
$ws.Range("L54").Value = 41999.001
$ws.Range("N54").Value = -43117.001

$ws = $wb.Worksheets.Item("GSM")
# Row 35
$ws.Range("H35").Value = 25000
$ws.Range("J35").Value = 25000
$ws.Range("L35").Value = 25000
$ws.Range("N35").Value = -25596

# Row 102
$ws.Range("H102").Value = 5504.524
$ws.Range("I102").Value = 6167.636
$ws.Range("K102").Value = 6167.636
$ws.Range("M102").Value = -4545.636

# Row 122
$ws.Range("H122").Value = 15883.2
$ws.Range("I122").Value = 23138
$ws.Range("J122").Value = 14069.5
$ws.Range("K122").Value = 69414
$ws.Range("L122").Value = 42208.5
$ws.Range("M122").Value = -66964
$ws.Range("N122").Value = -47108.5

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 19729.771
$ws.Range("I7").Value = 31585
$ws.Range("K7").Value = 31585
$ws.Range("M7").Value = -31473

# Row 22
$ws.Range("H22").Value = 22696.834
$ws.Range("J22").Value = 6952
$ws.Range("L22").Value = 6952
$ws.Range("N22").Value = -7542

# Row 27
$ws.Range("H27").Value = 22696.834
$ws.Range("J27").Value = 6952
$ws.Range("L27").Value = 6952
$ws.Range("N27").Value = -7166

# Row 68
$ws.Range("H68").Value = 4866.5
$ws.Range("I68").Value = 2999.5
$ws.Range("J68").Value = 5800
$ws.Range("K68").Value = 2999.5
$ws.Range("L68").Value = 5800
$ws.Range("M68").Value = -2250.5
$ws.Range("N68").Value = -7298

# Row 71
$ws.Range("H71").Value = 4866.5
$ws.Range("I71").Value = 2999.5
$ws.Range("J71").Value = 5800
$ws.Range("K71").Value = 14997.5
$ws.Range("L71").Value = 29000
$ws.Range("M71").Value = -11253.5
$ws.Range("N71").Value = -36488

# Row 100
$ws.Range("H100").Value = 2174.7
$ws.Range("I100").Value = 2174.7
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 2174.7
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -1633.7
$ws.Range("N100").ClearContents()

# Row 122
$ws.Range("H122").Value = 7575.9023
$ws.Range("I122").Value = 5905.9443
$ws.Range("K122").Value = 17717.8329
$ws.Range("M122").Value = -15267.8329

# Row 126
$ws.Range("H126").Value = 19729.771
$ws.Range("I126").Value = 31585
$ws.Range("K126").Value = 94755
$ws.Range("M126").Value = -92285

# Row 132
$ws.Range("H132").Value = 998529.4
$ws.Range("I132").Value = 2129564.5
$ws.Range("J132").Value = 8873.75
$ws.Range("K132").Value = 6388693.5
$ws.Range("L132").Value = 26621.25
$ws.Range("M132").Value = -6386163.5
$ws.Range("N132").Value = -31681.25

# Row 133
$ws.Range("H133").Value = 59999
$ws.Range("I133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("M133").ClearContents()

# Row 136
$ws.Range("H136").Value = 5599.75
$ws.Range("I136").Value = 2111.88
$ws.Range("K136").Value = 6335.64
$ws.Range("M136").Value = -3785.64
